$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change B2 from text "3" to a real number 3
$ws.Range("B2").Value = 3

# Add new row 3 data
$ws.Range("A3").Value = "Sunsi Wu"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "4"
$ws.Range("B3").ClearFormats()
$ws.Range("C3").Value = "should be"
$ws.Range("D3").Value = "DFT"
$ws.Range("E3").Value = "WRI"
$ws.Range("F3").Value = "c8048836-24fe-4e27-95aa-c7cfb58ac155"
$ws.Range("G3").Value = "rkc_hGb0Z_annotated.xlsx"
$ws.Range("H3").Value = "The structure of the global policies used in the experiments should be mentioned somewhere."
